# Slide 5, Shape 165 ("Shape 165" / body placeholder idx=1): the first
# paragraph currently reads "Pairplots" + ":" as two separate runs (the
# first run is flagged err="1" from a spell-check squiggle). The edit
# merges them into a single run reading "Pair plots / Scatter plots:"
# while keeping the second run's (non-err) formatting.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)

# Drop the first run ("Pairplots", chars 1-9) entirely so only the
# second run (":") survives, carrying its original formatting forward.
$oldRun = $para1.Characters(1, 9)
$oldRun.Text = ""

# Replace the remaining run's text (":") with the full desired text;
# since it is now the only run in the paragraph, this rewrites its
# <a:t> in place without creating any extra runs.
$remaining = $para1.Characters(1, 1)
$remaining.Text = "Pair plots / Scatter plots:"
